# Quarterly indexing esoteric bug-fix operation
# Updates the forecast error table (rows 2-11, columns B-G) on the
# active worksheet with corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.3349892545614097, 1.605491820351044, 4.859283619025187, 2.20437828401234, 2.205186593348543, 42)
    3  = @(-0.04710882026988346, 1.657196964163514, 4.932318376681804, 2.220882341926696, 2.247966106089669, 41)
    4  = @(0.4020606919227649, 1.506888104533207, 4.559370358069133, 2.135268216891998, 2.123788960030735, 40)
    5  = @(0.1019088716869291, 1.510897065334748, 4.536451325651763, 2.129894674778958, 2.155266318831541, 39)
    6  = @(0.3653246424028117, 1.548561036450066, 4.792308178577767, 2.189134116169626, 2.187409619450435, 38)
    7  = @(0.1349934622835307, 1.582982358002176, 4.917780729374886, 2.217606982622233, 2.244026807729286, 37)
    8  = @(0.3870259576278534, 1.497747264846744, 4.326038897300351, 2.079913194655092, 2.072575988579108, 36)
    9  = @(0.02232700200142348, 1.601898111008815, 4.9523141292121, 2.225379547226068, 2.257754918163958, 35)
    10 = @(0.4526287886684112, 1.583235315248334, 4.924710385845067, 2.219168850233138, 2.205190002309727, 34)
    11 = @(0.06611821871031362, 1.645663897126419, 4.676326697265005, 2.162481606225821, 2.194983752700733, 33)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
